# Decrease the "剩余" (remaining) value in column E by 1 for every data row
# (rows 2-99), except row 36 which is left unchanged in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value2 = $cell.Value2 - 1
}
